$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 69557
$ws.Range("E9").Value = 191325963

$ws.Range("C10").Value = 278095
$ws.Range("E10").Value = 1751062814

$ws.Range("C14").Value = 119177
$ws.Range("D14").Value = 25291
$ws.Range("E14").Value = 379185654

$ws.Range("C19").Value = 108882
$ws.Range("E19").Value = 344371442

$ws.Range("C46").Value = 55745
$ws.Range("E46").Value = 174066078

$ws.Range("C64").Value = 10535
$ws.Range("E64").Value = 25241365

$ws.Range("C81").Value = 26153
$ws.Range("E81").Value = 164984477

$ws.Range("C99").Value = 136530
$ws.Range("E99").Value = 862704807

$ws.Range("C102").Value = 13704
$ws.Range("E102").Value = 28900107

$ws.Range("C104").Value = 22087
$ws.Range("E104").Value = 84757147

$ws.Range("C108").Value = 45969
$ws.Range("E108").Value = 145891407

$ws.Range("C111").Value = 5909
$ws.Range("E111").Value = 11624445

$ws.Range("C115").Value = 17128
$ws.Range("E115").Value = 37534002

$ws.Range("C150").Value = 94996
$ws.Range("D150").Value = 21156
$ws.Range("E150").Value = 278333598

$ws.Range("C152").Value = 126011
$ws.Range("E152").Value = 715366771

$ws.Range("C156").Value = 47577
$ws.Range("E156").Value = 142329064

$ws.Range("C168").Value = 284699
$ws.Range("E168").Value = 1201325120

$ws.Range("C169").Value = 562449
$ws.Range("E169").Value = 1283421431

$ws.Range("C170").Value = 366857
$ws.Range("E170").Value = 2838343332

$ws.Range("C171").Value = 114999
$ws.Range("E171").Value = 441783336

$ws.Range("C174").Value = 356936
$ws.Range("E174").Value = 1012747979

$ws.Range("C175").Value = 125354
$ws.Range("E175").Value = 804226999

$ws.Range("C177").Value = 96726
$ws.Range("E177").Value = 174206346

$ws.Range("C179").Value = 235435
$ws.Range("E179").Value = 808022277

$ws.Range("C180").Value = 141376
$ws.Range("E180").Value = 338892343

$ws.Range("C182").Value = 6444
$ws.Range("E182").Value = 12761396

$ws.Range("C199").Value = 4083
$ws.Range("E199").Value = 8781595

$ws.Range("C203").Value = 12776
$ws.Range("E203").Value = 32128249

$ws.Range("C204").Value = 4642
$ws.Range("E204").Value = 11054908

$ws.Range("C205").Value = 10736
$ws.Range("E205").Value = 41639457

$ws.Range("C208").Value = 1516
$ws.Range("E208").Value = 3189944

$ws.Range("C213").Value = 3507
$ws.Range("E213").Value = 10616906

$ws.Range("C214").Value = 6102
$ws.Range("E214").Value = 10897376

$ws.Range("C247").Value = 29418
$ws.Range("E247").Value = 99422101

$ws.Range("C276").Value = 216591
$ws.Range("E276").Value = 1209582781

$ws.Range("C295").Value = 91314
$ws.Range("E295").Value = 552789957

$ws.Range("C298").Value = 11914
$ws.Range("E298").Value = 24041275

$ws.Range("C311").Value = 190831
$ws.Range("E311").Value = 585870572

$ws.Range("C313").Value = 220581
$ws.Range("E313").Value = 1369869685

$ws.Range("C320").Value = 67235
$ws.Range("E320").Value = 124544678

$ws.Range("C322").Value = 81132
$ws.Range("E322").Value = 254172109

$ws.Range("C323").Value = 94714
$ws.Range("E323").Value = 178746675
